# "actualizacion de febrero hay un archivo mal"
#
# The SIPOT "Reporte de Formatos" sheet was carrying 3er Trimestre 2021
# dates; bump the reporting period to 4to Trimestre 2021 (row 8) and
# leave the sheet scrolled/selected where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# Periodo que se informa: fecha de inicio / fecha de término (row 8)
# Stored as plain date-serial numbers so the existing m/d/yyyy cell
# format (style index 7) renders them correctly.
$ws.Range("B8").Value = 44470   # 10/1/2021  (was 7/1/2021)
$ws.Range("C8").Value = 44561   # 12/31/2021 (was 9/30/2021)

# Fecha de validación / Fecha de actualización (row 8)
$ws.Range("S8").Value = 44571   # 1/10/2022 (was 10/11/2021)
$ws.Range("T8").Value = 44571   # 1/10/2022 (was 10/11/2021)

# Leave the window scrolled over to the right-hand columns and the
# cursor parked on T13, matching where the author's session was left.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 17
$ws.Range("T13").Select()
